$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.036645423468762
$ws.Cells.Item(2, 4).Value = 1.036688405373256
$ws.Cells.Item(2, 5).Value = 1.040294338750955
$ws.Cells.Item(2, 6).Value = 1.035334893408444
$ws.Cells.Item(2, 9).Value = 1.033236240024902
$ws.Cells.Item(2, 10).Value = 1.041752610746729
$ws.Cells.Item(2, 11).Value = 1.03948143130787
$ws.Cells.Item(2, 12).Value = 1.043077095115534
$ws.Cells.Item(2, 13).Value = 1.038131794151826
$ws.Cells.Item(2, 14).Value = 1.043232018944731

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.039020264665833
$ws.Cells.Item(3, 4).Value = 1.038463267819327
$ws.Cells.Item(3, 5).Value = 1.042620401078498
$ws.Cells.Item(3, 6).Value = 1.038359239604055
$ws.Cells.Item(3, 9).Value = 1.033891104299029
$ws.Cells.Item(3, 10).Value = 1.043764338209669
$ws.Cells.Item(3, 11).Value = 1.041063389808965
$ws.Cells.Item(3, 12).Value = 1.045209568271824
$ws.Cells.Item(3, 13).Value = 1.040959636942425
$ws.Cells.Item(3, 14).Value = 1.045246603291369

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.040545882555764
$ws.Cells.Item(4, 4).Value = 1.039601999471175
$ws.Cells.Item(4, 5).Value = 1.044114900361303
$ws.Cells.Item(4, 6).Value = 1.040303805237578
$ws.Cells.Item(4, 9).Value = 1.034308750621779
$ws.Cells.Item(4, 10).Value = 1.045055070990584
$ws.Cells.Item(4, 11).Value = 1.042076889080467
$ws.Cells.Item(4, 12).Value = 1.046578515052757
$ws.Cells.Item(4, 13).Value = 1.042776934555341
$ws.Cells.Item(4, 14).Value = 1.046539169060883

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.041184667201573
$ws.Cells.Item(5, 4).Value = 1.040078442291129
$ws.Cells.Item(5, 5).Value = 1.044740706839957
$ws.Cells.Item(5, 6).Value = 1.041118419048466
$ws.Cells.Item(5, 9).Value = 1.034482891271994
$ws.Cells.Item(5, 10).Value = 1.045595120305858
$ws.Cells.Item(5, 11).Value = 1.042500583867387
$ws.Cells.Item(5, 12).Value = 1.047151467530324
$ws.Cells.Item(5, 13).Value = 1.043538012476241
$ws.Cells.Item(5, 14).Value = 1.047079985308129

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.041291772185177
$ws.Cells.Item(6, 4).Value = 1.040158306810932
$ws.Cells.Item(6, 5).Value = 1.044845638769174
$ws.Cells.Item(6, 6).Value = 1.041255029736353
$ws.Cells.Item(6, 9).Value = 1.034512046545569
$ws.Cells.Item(6, 10).Value = 1.045685647591342
$ws.Cells.Item(6, 11).Value = 1.042571585817721
$ws.Cells.Item(6, 12).Value = 1.047247520749457
$ws.Cells.Item(6, 13).Value = 1.043665632226443
$ws.Cells.Item(6, 14).Value = 1.047170641152742

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.040554428101116
$ws.Cells.Item(7, 4).Value = 1.039608374618814
$ws.Cells.Item(7, 5).Value = 1.044123272085813
$ws.Cells.Item(7, 6).Value = 1.040314701358298
$ws.Cells.Item(7, 9).Value = 1.034311083119337
$ws.Cells.Item(7, 10).Value = 1.045062297198524
$ws.Cells.Item(7, 11).Value = 1.04208255980329
$ws.Cells.Item(7, 12).Value = 1.046586180825256
$ws.Cells.Item(7, 13).Value = 1.042787115459999
$ws.Cells.Item(7, 14).Value = 1.046546405530866

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.037450341448225
$ws.Cells.Item(8, 4).Value = 1.037290270707987
$ws.Cells.Item(8, 5).Value = 1.041082678206691
$ws.Cells.Item(8, 6).Value = 1.036359601760893
$ws.Cells.Item(8, 9).Value = 1.033458829316469
$ws.Cells.Item(8, 10).Value = 1.042434793038949
$ws.Cells.Item(8, 11).Value = 1.040018187638834
$ws.Cells.Item(8, 12).Value = 1.043800067923118
$ws.Cells.Item(8, 13).Value = 1.039090114869248
$ws.Cells.Item(8, 14).Value = 1.043915170014007

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.031892914208762
$ws.Cells.Item(9, 4).Value = 1.033128825304128
$ws.Cells.Item(9, 5).Value = 1.035640645218081
$ws.Cells.Item(9, 6).Value = 1.029291476428755
$ws.Cells.Item(9, 9).Value = 1.031909400600137
$ws.Cells.Item(9, 10).Value = 1.037718062338713
$ws.Cells.Item(9, 11).Value = 1.036300797757999
$ws.Cells.Item(9, 12).Value = 1.038804395703065
$ws.Cells.Item(9, 13).Value = 1.032476094703013
$ws.Cells.Item(9, 14).Value = 1.03919174101516

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.028124919526673
$ws.Cells.Item(10, 4).Value = 1.030299884804251
$ws.Cells.Item(10, 5).Value = 1.031952093247566
$ws.Cells.Item(10, 6).Value = 1.02450752621088
$ws.Cells.Item(10, 9).Value = 1.030843027905995
$ws.Cells.Item(10, 10).Value = 1.034511607118213
$ws.Cells.Item(10, 11).Value = 1.033765973609437
$ws.Cells.Item(10, 12).Value = 1.035412217133781
$ws.Cells.Item(10, 13).Value = 1.027994694907158
$ws.Cells.Item(10, 14).Value = 1.03598073226046

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.026477443597654
$ws.Cells.Item(11, 4).Value = 1.029061253280891
$ws.Cells.Item(11, 5).Value = 1.030339647589795
$ws.Cells.Item(11, 6).Value = 1.022417728059009
$ws.Cells.Item(11, 9).Value = 1.030373034026487
$ws.Cells.Item(11, 10).Value = 1.033107646084166
$ws.Cells.Item(11, 11).Value = 1.032654271977982
$ws.Cells.Item(11, 12).Value = 1.033927862639007
$ws.Cells.Item(11, 13).Value = 1.02603592463209
$ws.Cells.Item(11, 14).Value = 1.034574777440709

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.025863025174161
$ws.Cells.Item(12, 4).Value = 1.028599052780142
$ws.Cells.Item(12, 5).Value = 1.029738338805551
$ws.Cells.Item(12, 6).Value = 1.021638622784307
$ws.Cells.Item(12, 9).Value = 1.030197190534821
$ws.Cells.Item(12, 10).Value = 1.032583744057222
$ws.Cells.Item(12, 11).Value = 1.032239157721607
$ws.Cells.Item(12, 12).Value = 1.03337410022579
$ws.Cells.Item(12, 13).Value = 1.025305497808029
$ws.Cells.Item(12, 14).Value = 1.034050131412795

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.025994933244908
$ws.Cells.Item(13, 4).Value = 1.028698293198765
$ws.Cells.Item(13, 5).Value = 1.029867430337293
$ws.Cells.Item(13, 6).Value = 1.02180587477639
$ws.Cells.Item(13, 9).Value = 1.030234967342141
$ws.Cells.Item(13, 10).Value = 1.032696232999175
$ws.Cells.Item(13, 11).Value = 1.032328300706567
$ws.Cells.Item(13, 12).Value = 1.033492994304166
$ws.Cells.Item(13, 13).Value = 1.025462307614007
$ws.Cells.Item(13, 14).Value = 1.034162780101947

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.02642670650191
$ws.Cells.Item(14, 4).Value = 1.029023091207208
$ws.Cells.Item(14, 5).Value = 1.030289992127255
$ws.Cells.Item(14, 6).Value = 1.022353385985458
$ws.Cells.Item(14, 9).Value = 1.030358524755108
$ws.Cells.Item(14, 10).Value = 1.033064389738978
$ws.Cells.Item(14, 11).Value = 1.032620003300524
$ws.Cells.Item(14, 12).Value = 1.033882138013108
$ws.Cells.Item(14, 13).Value = 1.025975606120523
$ws.Cells.Item(14, 14).Value = 1.03453145966655

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.026692406124531
$ws.Cells.Item(15, 4).Value = 1.029222927490947
$ws.Cells.Item(15, 5).Value = 1.030550029317333
$ws.Cells.Item(15, 6).Value = 1.022690343251967
$ws.Cells.Item(15, 9).Value = 1.030434483843625
$ws.Cells.Item(15, 10).Value = 1.03329090190617
$ws.Cells.Item(15, 11).Value = 1.032799440380375
$ws.Cells.Item(15, 12).Value = 1.034121581041938
$ws.Cells.Item(15, 13).Value = 1.026291485157247
$ws.Cells.Item(15, 14).Value = 1.034758293507

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.0282339149109
$ws.Cells.Item(16, 4).Value = 1.030381795180968
$ws.Cells.Item(16, 5).Value = 1.032058777340229
$ws.Cells.Item(16, 6).Value = 1.024645824061635
$ws.Cells.Item(16, 9).Value = 1.030874043735988
$ws.Cells.Item(16, 10).Value = 1.034604449694778
$ws.Cells.Item(16, 11).Value = 1.033839451345219
$ws.Cells.Item(16, 12).Value = 1.035510395488268
$ws.Cells.Item(16, 13).Value = 1.028124297799901
$ws.Cells.Item(16, 14).Value = 1.036073706684134

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.029196545840754
$ws.Cells.Item(17, 4).Value = 1.031105015993786
$ws.Cells.Item(17, 5).Value = 1.033001029371601
$ws.Cells.Item(17, 6).Value = 1.025867465461939
$ws.Cells.Item(17, 9).Value = 1.031147540610749
$ws.Cells.Item(17, 10).Value = 1.035424190849806
$ws.Cells.Item(17, 11).Value = 1.034488004656033
$ws.Cells.Item(17, 12).Value = 1.036377354583811
$ws.Cells.Item(17, 13).Value = 1.029269002391709
$ws.Cells.Item(17, 14).Value = 1.036894611965627

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.029756502819998
$ws.Cells.Item(18, 4).Value = 1.031525542656419
$ws.Cells.Item(18, 5).Value = 1.033549160363076
$ws.Cells.Item(18, 6).Value = 1.026578268516091
$ws.Cells.Item(18, 9).Value = 1.031306272861843
$ws.Cells.Item(18, 10).Value = 1.035900837176976
$ws.Cells.Item(18, 11).Value = 1.034864937898271
$ws.Cells.Item(18, 12).Value = 1.036881545020283
$ws.Cells.Item(18, 13).Value = 1.0299349311986
$ws.Cells.Item(18, 14).Value = 1.037371935185254

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.029947176735989
$ws.Cells.Item(19, 4).Value = 1.031668710055397
$ws.Cells.Item(19, 5).Value = 1.033735812124473
$ws.Cells.Item(19, 6).Value = 1.026820338783601
$ws.Cells.Item(19, 9).Value = 1.03136026263661
$ws.Cells.Item(19, 10).Value = 1.036063109891777
$ws.Cells.Item(19, 11).Value = 1.03499323401338
$ws.Cells.Item(19, 12).Value = 1.037053210096887
$ws.Cells.Item(19, 13).Value = 1.030161700626133
$ws.Cells.Item(19, 14).Value = 1.037534438345923

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.029093423429009
$ws.Cells.Item(20, 4).Value = 1.031027557788461
$ws.Cells.Item(20, 5).Value = 1.032900087124858
$ws.Cells.Item(20, 6).Value = 1.025736577697579
$ws.Cells.Item(20, 9).Value = 1.031118279294255
$ws.Cells.Item(20, 10).Value = 1.035336395519551
$ws.Cells.Item(20, 11).Value = 1.034418561802679
$ws.Cells.Item(20, 12).Value = 1.036284492956282
$ws.Cells.Item(20, 13).Value = 1.029146369003636
$ws.Cells.Item(20, 14).Value = 1.036806691955933

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.026299628964464
$ws.Cells.Item(21, 4).Value = 1.028927505252086
$ws.Cells.Item(21, 5).Value = 1.030165624411772
$ws.Cells.Item(21, 6).Value = 1.022192237456771
$ws.Cells.Item(21, 9).Value = 1.030322175374072
$ws.Cells.Item(21, 10).Value = 1.032956043833949
$ws.Cells.Item(21, 11).Value = 1.032534164787485
$ws.Cells.Item(21, 12).Value = 1.03376761195321
$ws.Cells.Item(21, 13).Value = 1.025824532089398
$ws.Cells.Item(21, 14).Value = 1.034422959897908

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.024528711956718
$ws.Cells.Item(22, 4).Value = 1.02759483510967
$ws.Cells.Item(22, 5).Value = 1.02843257924392
$ws.Cells.Item(22, 6).Value = 1.01994715973884
$ws.Cells.Item(22, 9).Value = 1.02981428940944
$ws.Cells.Item(22, 10).Value = 1.031445450503807
$ws.Cells.Item(22, 11).Value = 1.031336733257657
$ws.Cells.Item(22, 12).Value = 1.032171183561837
$ws.Cells.Item(22, 13).Value = 1.023719405223801
$ws.Cells.Item(22, 14).Value = 1.032910221351969

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.025468896609652
$ws.Cells.Item(23, 4).Value = 1.028302494295541
$ws.Cells.Item(23, 5).Value = 1.029352632482194
$ws.Cells.Item(23, 6).Value = 1.021138929694556
$ws.Cells.Item(23, 9).Value = 1.03008423494428
$ws.Cells.Item(23, 10).Value = 1.032247593845893
$ws.Cells.Item(23, 11).Value = 1.031972732432392
$ws.Cells.Item(23, 12).Value = 1.033018829970796
$ws.Cells.Item(23, 13).Value = 1.024836977529632
$ws.Cells.Item(23, 14).Value = 1.033713503829609

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.029140024729226
$ws.Cells.Item(24, 4).Value = 1.031062561879764
$ws.Cells.Item(24, 5).Value = 1.032945703115043
$ws.Cells.Item(24, 6).Value = 1.025795725675636
$ws.Cells.Item(24, 9).Value = 1.031131503675116
$ws.Cells.Item(24, 10).Value = 1.035376071065173
$ws.Cells.Item(24, 11).Value = 1.034449944228636
$ws.Cells.Item(24, 12).Value = 1.036326457735144
$ws.Cells.Item(24, 13).Value = 1.029201787179448
$ws.Cells.Item(24, 14).Value = 1.036846423845381

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.033340445323225
$ws.Cells.Item(25, 4).Value = 1.034214057275931
$ws.Cells.Item(25, 5).Value = 1.037057912563369
$ws.Cells.Item(25, 6).Value = 1.031131012683623
$ws.Cells.Item(25, 9).Value = 1.038948099021789
$ws.Cells.Item(25, 10).Value = 1.038948099021789
$ws.Cells.Item(25, 11).Value = 1.037271579082289
$ws.Cells.Item(25, 12).Value = 1.040106493825893
$ws.Cells.Item(25, 13).Value = 1.034198287926408
$ws.Cells.Item(25, 14).Value = 1.040423524491412
